# Update rss_feed worksheet: row 94/95 "File Name" become real numbers,
# and append the two new RSS items as rows 96/97.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 94 / 95: column A ("File Name") switches from text "9" to the number 9.
$ws.Range("A94").Value = 9
$ws.Range("A95").Value = 9

# ---- New row 96 ----
$ws.Range("A96").NumberFormat = "@"
$ws.Range("A96").Value = '11'
$ws.Range("B96").Value = 'Marketing Video for Land Parcel - Upwork'
$ws.Range("C96").Value = 'https://www.upwork.com/jobs/Marketing-Video-for-Land-Parcel_%7E017bee3c8c9dc0a1ef?source=rss'
$ws.Range("F96").Value = 'Sat, 15 Jun 2024 03:00:57 +0000'
$ws.Range("G96").Value = 'https://www.upwork.com/jobs/Marketing-Video-for-Land-Parcel_%7E017bee3c8c9dc0a1ef?source=rss'
$ws.Range("J96").Value = 'June 15, 2024 03:00 UTC'
$ws.Range("K96").Value = 'Video Production'
$ws.Range("L96").Value = 'Nonfiction,     Video Commercial,     Real Estate Video'
$ws.Range("M96").Value = 'Mexico'
$ws.Range("D96").Value = @'
Need a short 1-2 min video that features a piece of land that is for sale. the video should highlight the strategic benefits that come with owning the land such as location, area demographics, nearby cities and ports. It should highlight the surface of the land and also talk a little about the project that is being developed on the land. 
We are project developers of an industrial park and are selling a part of the project while the land is still, with no infrastructure to fund the development process.
The whole project is the 400+ hectare plot in the KMZ and the land that will be sold is labeled &quot;FINCA 24744&quot; 
Budget
: $150
Posted On
: June 15, 2024 03:00 UTC
Category
: Video Production
Skills
:Nonfiction,     Video Commercial,     Real Estate Video    
Skills
:        Nonfiction,                     Video Commercial,                     Real Estate Video            
Country
: Mexico
click to apply

'@
$ws.Range("E96").Value = @'
Need a short 1-2 min video that features a piece of land that is for sale. the video should highlight the strategic benefits that come with owning the land such as location, area demographics, nearby cities and ports. It should highlight the surface of the land and also talk a little about the project that is being developed on the land. <br /><br />
We are project developers of an industrial park and are selling a part of the project while the land is still, with no infrastructure to fund the development process.<br /><br />
The whole project is the 400+ hectare plot in the KMZ and the land that will be sold is labeled &amp;quot;FINCA 24744&amp;quot; <br /><br /><b>Budget</b>: $150
<br /><b>Posted On</b>: June 15, 2024 03:00 UTC<br /><b>Category</b>: Video Production<br /><b>Skills</b>:Nonfiction,     Video Commercial,     Real Estate Video    
<br /><b>Skills</b>:        Nonfiction,                     Video Commercial,                     Real Estate Video            <br /><b>Country</b>: Mexico
<br /><a href="https://www.upwork.com/jobs/Marketing-Video-for-Land-Parcel_%7E017bee3c8c9dc0a1ef?source=rss">click to apply</a>

'@
$ws.Range("H96").NumberFormat = "@"
$ws.Range("H96").Value = ""
$ws.Range("I96").NumberFormat = "@"
$ws.Range("I96").Value = '$150'
$ws.Rows(96).EntireRow.AutoFit()

# ---- New row 97 ----
$ws.Range("A97").NumberFormat = "@"
$ws.Range("A97").Value = '11'
$ws.Range("B97").Value = 'Join 2 short videos into one - Upwork'
$ws.Range("C97").Value = 'https://www.upwork.com/jobs/Join-short-videos-into-one_%7E01b5a1dba74a8955db?source=rss'
$ws.Range("F97").Value = 'Sat, 15 Jun 2024 02:35:52 +0000'
$ws.Range("G97").Value = 'https://www.upwork.com/jobs/Join-short-videos-into-one_%7E01b5a1dba74a8955db?source=rss'
$ws.Range("J97").Value = 'June 15, 2024 02:35 UTC'
$ws.Range("K97").Value = 'Video Editing'
$ws.Range("L97").Value = 'Video Editing,     Adobe Premiere Pro,     Video Post-Editing'
$ws.Range("M97").Value = 'United States'
$ws.Range("D97").Value = @'
I need to create a video for a presentation that involves the assembly and final walk through of a folding container house.  This video is just for a presentation to my co-workers and will not be seen by others.  I have a 10-15 second video taken with my I-Phone of the INTERIOR of the completed house.  However, we did not get any video of the assembly of the house when it was delivered, so I found a generic video (artists rendition) that shows what I want, so I'd like that added to MY interior video.  The who length of the video will be less than 30 seconds.  I will need you to join the two together and then add some written text on the screen.  The attached PHOTO is a screenshot of the youtube video I want to use.
Can this be done rather easily?
Budget
: $50
Posted On
: June 15, 2024 02:35 UTC
Category
: Video Editing
Skills
:Video Editing,     Adobe Premiere Pro,     Video Post-Editing    
Skills
:        Video Editing,                     Adobe Premiere Pro,                     Video Post-Editing            
Country
: United States
click to apply

'@
$ws.Range("E97").Value = @'
I need to create a video for a presentation that involves the assembly and final walk through of a folding container house.&nbsp;&nbsp;This video is just for a presentation to my co-workers and will not be seen by others.&nbsp;&nbsp;I have a 10-15 second video taken with my I-Phone of the INTERIOR of the completed house.&nbsp;&nbsp;However, we did not get any video of the assembly of the house when it was delivered, so I found a generic video (artists rendition) that shows what I want, so I&#039;d like that added to MY interior video.&nbsp;&nbsp;The who length of the video will be less than 30 seconds.&nbsp;&nbsp;I will need you to join the two together and then add some written text on the screen.&nbsp;&nbsp;The attached PHOTO is a screenshot of the youtube video I want to use.<br />
Can this be done rather easily?<br /><br /><b>Budget</b>: $50
<br /><b>Posted On</b>: June 15, 2024 02:35 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing,     Adobe Premiere Pro,     Video Post-Editing    
<br /><b>Skills</b>:        Video Editing,                     Adobe Premiere Pro,                     Video Post-Editing            <br /><b>Country</b>: United States
<br /><a href="https://www.upwork.com/jobs/Join-short-videos-into-one_%7E01b5a1dba74a8955db?source=rss">click to apply</a>

'@
$ws.Range("H97").NumberFormat = "@"
$ws.Range("H97").Value = ""
$ws.Range("I97").NumberFormat = "@"
$ws.Range("I97").Value = '$50'
$ws.Rows(97).EntireRow.AutoFit()

